# "new csvs April 25" -- append the latest basin record (BENT / Benito-Ntem)
# that was exported from the April 25 CSV refresh, and leave the cursor
# where the user last clicked before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "BENT"
$ws.Range("C4").Value = "Benito/Ntem"
$ws.Range("D4").Value = (Get-Date -Year 2023 -Month 7 -Day 2).Date
$ws.Range("E4").Value = (Get-Date -Year 2024 -Month 3 -Day 15).Date

$ws.Range("D12").Select()
